$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "Menu bar shows ... User" run with the following " - " run
#    into a single run (Find/Replace with identical text normalises the two
#    adjacent runs that share the same formatting into one run).
# ---------------------------------------------------------------------------
$menuText = "Menu bar shows extra options being available to the User " + [char]0x2013 + " "
$null = $d.Content.Find.Execute($menuText, $true, $false, $false, $false, $false, $true, 1, $false, $menuText, 2)

# ---------------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark (it sat at the end of the bullet-point
#    paragraph just before the section break).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) Fix the password value "Super123" -> "super123" (lower-case the first
#    letter). Re-typing that single character is what leaves the "_GoBack"
#    bookmark positioned right after it, splitting the run into "s" / "uper123".
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Super123")
if ($found) {
    $firstChar = $d.Range($r.Start, $r.Start + 1)
    $firstChar.Text = "s"

    $bmPos = $r.Start + 1
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 4) Update the cached NUMPAGES field result in the default footer from "2"
#    to "1".
# ---------------------------------------------------------------------------
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$fields = $ftr.Range.Fields
for ($i = 1; $i -le $fields.Count; $i++) {
    $fld = $fields.Item($i)
    if ($fld.Code.Text -match "NUMPAGES") {
        $fres = $fld.Result
        $null = $fres.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)
    }
}
